$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Total time taken for the ride"
$ws.Cells.Item(1, 2).Value = 0.03498324074074075

$ws.Cells.Item(2, 1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(2, 2).Value = 33.50112111111111

$ws.Cells.Item(3, 1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(3, 2).Value = 1707.780862160555

$ws.Cells.Item(4, 1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(4, 2).Value = 39.51

$ws.Cells.Item(5, 1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(5, 2).Value = 6.548

$ws.Cells.Item(6, 1).Value = "Starting SoC (%)"
$ws.Cells.Item(6, 2).Value = 99

$ws.Cells.Item(7, 1).Value = "Ending SoC (%)"
$ws.Cells.Item(7, 2).Value = 16

$ws.Cells.Item(8, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(8, 2).Value = 32.7550019235124

$ws.Cells.Item(9, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(9, 2).Value = 52.13801745908815

$ws.Cells.Item(10, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(10, 2).Value = 83

$ws.Cells.Item(11, 1).Value = "Mode"
$ws.Cells.Item(11, 2).Value = "Custom mode`n78.47%`nEco mode`n21.53%"

$ws.Cells.Item(12, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(12, 2).Value = 5428.3737

$ws.Cells.Item(13, 1).Value = "Average Power(kW)"
$ws.Cells.Item(13, 2).Value = -2041.172345211819

$ws.Cells.Item(14, 1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(14, 2).Value = 0.06724817500000001

$ws.Cells.Item(15, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(15, 2).Value = 0.003937596943956988

$ws.Cells.Item(16, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(16, 2).Value = 3.43

$ws.Cells.Item(17, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.088

$ws.Cells.Item(18, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(18, 2).Value = 0.3420000000000001

$ws.Cells.Item(19, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(19, 2).Value = 38

$ws.Cells.Item(20, 1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(20, 2).Value = 47

$ws.Cells.Item(21, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(21, 2).Value = 9

$ws.Cells.Item(22, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(22, 2).Value = 68

$ws.Cells.Item(23, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(23, 2).Value = 67

$ws.Cells.Item(24, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(24, 2).Value = 65

$ws.Cells.Item(25, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(25, 2).Value = 71

$ws.Cells.Item(26, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(26, 2).Value = 101

$ws.Cells.Item(27, 1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(27, 2).Value = 0

$ws.Cells.Item(28, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(28, 2).Value = 47

$ws.Cells.Item(29, 1).Value = "lowest cell temp(C)"
$ws.Cells.Item(29, 2).Value = 38

$ws.Cells.Item(30, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(30, 2).Value = 9

$ws.Cells.Item(31, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(31, 2).Value = 54

$ws.Cells.Item(32, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(32, 2).Value = 1.80906054

$ws.Cells.Item(33, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(33, 2).Value = 0.000000166286173615707

$ws.Cells.Item(34, 1).Value = "Idling time percentage"
$ws.Cells.Item(34, 2).Value = 0.6348974677078012

$ws.Cells.Item(35, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(35, 2).Value = 4.84930307231993

$ws.Cells.Item(36, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(36, 2).Value = 8.115011311391667

$ws.Cells.Item(37, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(37, 2).Value = 16.06947383784573

$ws.Cells.Item(38, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(38, 2).Value = 26.73137269211121

$ws.Cells.Item(39, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(39, 2).Value = 12.5702400934102

$ws.Cells.Item(40, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(40, 2).Value = 13.05553528424433

$ws.Cells.Item(41, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(41, 2).Value = 10.07808509085602

$ws.Cells.Item(42, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(42, 2).Value = 7.724585857111581

$ws.Cells.Item(43, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(43, 2).Value = 0
